# Update "想去人数" (interest count) figures for several events.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5486
$ws1.Range("F12").Value = 3001
$ws1.Range("F14").Value = 1605

# Sheet "全部类型" (All types) — same events duplicated here with shifted rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5486
$ws4.Range("F13").Value = 3001
$ws4.Range("F15").Value = 1605
